$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: ATROVENT 250MCG/2ML 20 UNIT DOSE VIAL ---
# H9 (balance) is plain text (numFmt "@") - direct text assignment keeps it as text.
$ws.Range("H9").Value = "3:14"
# P9 (selling price) has a numeric format ("0.00"), so a numeric-looking string would be
# auto-converted to a real number by Excel. Temporarily switch the format to Text so the
# value is stored as text (shared string), then restore the original numeric format so the
# cell style is unchanged.
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "14.3000"
$ws.Range("P9").NumberFormat = "0.00"
# Q9 (transaction count) is plain text.
$ws.Range("Q9").Value = "0:1"

# --- Row 19: PULMICORT 0.5MG/ML 20 NEBULIZER VIAL SUSP. ---
$ws.Range("H19").Value = "1:4"
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value = "75.2000"
$ws.Range("P19").NumberFormat = "0.00"
$ws.Range("Q19").Value = "0:2"

# --- Row 31: سرنجات 3 سم (Syringes 3cm) ---
$ws.Range("P31").NumberFormat = "@"
$ws.Range("P31").Value = "8.0000"
$ws.Range("P31").NumberFormat = "0.00"
$ws.Range("Q31").Value = "4:0"

# --- Row 39: recalculated total of the "selling price" column ---
$ws.Range("P39").Value = 988.77
